$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC")

$ws.Range("E4").Value = $ws.Range("D4").Text
$ws.Range("F4").Value = "Passed"
